# Updates cryptos list: refreshed prices/volume percentages for Wed Nov 27 2024 run,
# plus two rank swaps (EthereumClassic/PolygonEcosystemToken and OKB/Cosmos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '94.698.71'
$ws.Range('E2').Value = '  +1.91%  '

$ws.Range('D3').Value = '3.501.89'
$ws.Range('E3').Value = '  +4.82%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.31'
$ws.Range('E5').Value = '  +3.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '629.18'
$ws.Range('E6').Value = '  +1.61%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.45'
$ws.Range('E7').Value = '  +5.83%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.396'
$ws.Range('E8').Value = '  +3.20%  '

$ws.Range('E9').Value = '  +0.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.01'
$ws.Range('E10').Value = '  +8.52%  '

$ws.Range('D11').Value = '3.499.33'
$ws.Range('E11').Value = '  +4.64%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.65'
$ws.Range('E12').Value = '  +3.82%  '

$ws.Range('E13').Value = '  +5.03%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.28'
$ws.Range('E14').Value = '  +4.92%  '

$ws.Range('D15').Value = '4.162.46'
$ws.Range('E15').Value = '  +4.92%  '

$ws.Range('D16').Value = '94.522.05'
$ws.Range('E16').Value = '  +1.83%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000251'
$ws.Range('E17').Value = '  +3.78%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.39'
$ws.Range('E18').Value = '  +5.34%  '

$ws.Range('D19').Value = '3.502.32'
$ws.Range('E19').Value = '  +4.88%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.69'
$ws.Range('E20').Value = '  +13.86%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.01'
$ws.Range('E21').Value = '  +4.14%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.504'
$ws.Range('E22').Value = '  +12.27%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '520.86'
$ws.Range('E23').Value = '  +6.49%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.39'
$ws.Range('E24').Value = '  +2.63%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.73'
$ws.Range('E25').Value = '  +10.61%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000186'
$ws.Range('E26').Value = '  +2.93%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.42'
$ws.Range('E27').Value = '  +7.73%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.25'
$ws.Range('E28').Value = '  +5.87%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.95'
$ws.Range('E29').Value = '  +11.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.53'
$ws.Range('E30').Value = '  +3.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.141'
$ws.Range('E31').Value = '  +4.88%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  +0.02%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.182'
$ws.Range('E33').Value = '  +5.80%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.992'
$ws.Range('E34').Value = '  -0.27%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.84'
$ws.Range('E35').Value = '  +5.51%  '

$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.562'
$ws.Range('E36').Value = '  +6.57%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '583.25'
$ws.Range('E37').Value = '  +10.38%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.45'
$ws.Range('E38').Value = '  +6.62%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.58'
$ws.Range('E39').Value = '  +3.46%  '

$ws.Range('E41').Value = '  +5.02%  '

$ws.Range('E42').Value = '  +2.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0425'
$ws.Range('E43').Value = '  +5.70%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.75'
$ws.Range('E44').Value = '  -1.16%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.71'
$ws.Range('E45').Value = '  +1.86%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.55'
$ws.Range('E46').Value = '  +3.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.55'
$ws.Range('E47').Value = '  -0.14%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.17'
$ws.Range('E48').Value = '  +2.66%  '

$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.72'
$ws.Range('E49').Value = '  +3.13%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.22'
$ws.Range('E50').Value = '  +4.59%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.08'
$ws.Range('E51').Value = '  +1.15%  '
